# The commit adds one new daily-price record for Kiwi "Primera" quality at
# "Región de O'Higgins" (date 45131) at the top of this block of rows, pushing
# every existing record in rows 864-985 down by one row (to 865-986).
#
# We reproduce that by inserting a blank row at row 864 (which shifts
# rows 864:985 down to 865:986, automatically growing the sheet's used range
# to row 986) and then populating the newly-inserted row 864 with the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 864; everything from 864 downward moves to 865+.
$ws.Rows("864:864").Insert()

# Populate the new row 864 with the new record.
$ws.Range("A864").Value2 = 9
$ws.Range("B864").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C864").Value2 = "Metropolitana"
$ws.Range("D864").Value2 = 45131
$ws.Range("E864").Value2 = 13
$ws.Range("F864").Value2 = "Fruta"
$ws.Range("G864").Value2 = 100101
$ws.Range("H864").Value2 = "Berries"
$ws.Range("I864").Value2 = 100101007
$ws.Range("J864").Value2 = "Kiwi"
$ws.Range("K864").Value2 = "Hayward"
$ws.Range("L864").Value2 = "Primera"
$ws.Range("M864").Value2 = 280
$ws.Range("N864").Value2 = 5500
$ws.Range("O864").Value2 = 5500
$ws.Range("P864").Value2 = 5500
$ws.Range("Q864").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R864").Value2 = "Región de O'Higgins"
$ws.Range("S864").Value2 = 550
$ws.Range("T864").Value2 = 10
